# Weekly update: a new price reading is inserted at the top of the data
# block (row 89) and all subsequent rows shift down by one, with the
# previously-last row (203) re-appearing (copied in full) as the new
# last row (204).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 89
$lastRow  = 203
$newLastRow = 204

# Columns that vary from row to row (the "record" columns).
$varCols = @(4, 10, 11, 12, 13, 14, 15, 16, 17)   # D, J, K, L, M, N, O, P, Q

# 1) Capture the current ("before") values of the varying columns for
#    every row in the block, rows 89..203, keyed by row number.
$captured = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $rowVals = @{}
    foreach ($c in $varCols) {
        $rowVals[$c] = $ws.Cells.Item($r, $c).Value2
    }
    $captured[$r] = $rowVals
}

# 2) Create the new last row (204) as a full copy of the old last row
#    (203) - every column A..R. Only column D (the date column) carries
#    an explicit number format in this sheet, so only it needs its
#    format copied across; leave the rest with the default style.
for ($c = 1; $c -le 18; $c++) {
    $srcCell = $ws.Cells.Item($lastRow, $c)
    $dstCell = $ws.Cells.Item($newLastRow, $c)
    $dstCell.Value = $srcCell.Value2
    if ($c -eq 4) {
        $dstCell.NumberFormat = $srcCell.NumberFormat
    }
}

# 3) Shift every row down by one: row r (90..203) gets the varying
#    column values that used to live in row r-1 (captured above).
for ($r = $lastRow; $r -ge ($firstRow + 1); $r--) {
    $src = $captured[$r - 1]
    foreach ($c in $varCols) {
        $ws.Cells.Item($r, $c).Value = $src[$c]
    }
}

# 4) Row 89 receives the brand-new price reading.
$ws.Cells.Item($firstRow, 4).Value  = 44482   # D - Fecha
$ws.Cells.Item($firstRow, 10).Value = 60      # J - Volumen
$ws.Cells.Item($firstRow, 11).Value = 3000    # K - Precio minimo
$ws.Cells.Item($firstRow, 12).Value = 4000    # L - Precio maximo
$ws.Cells.Item($firstRow, 13).Value = 3333    # M - Precio promedio ponderado
$ws.Cells.Item($firstRow, 14).Value = "$/docena de atados (3 kilos)"  # N - Unidad
$ws.Cells.Item($firstRow, 15).Value = "Provincia de Cautín"            # O - Origen
$ws.Cells.Item($firstRow, 16).Value = 1111    # P - Precio $/Kg
$ws.Cells.Item($firstRow, 17).Value = 3       # Q - Kg o Unidades
